$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Insert two new rows above the current row 10 ("office\Users" / office_users),
# which (like the rest of column A in this block) carries style s="7" on A10,
# so the freshly inserted A10/A11 cells pick that style up automatically.
$ws.Range("A10:A11").EntireRow.Insert()

# New "admin_user_personal" / "admin_user_legals" sub-rows (bold style in column C,
# like the existing "Type (Comfort)" style cells elsewhere on the sheet).
$ws.Range("C10").Value = "admin_user_personal"
$ws.Range("C10").Font.Bold = $true
$ws.Range("C11").Value = "admin_user_legals"
$ws.Range("C11").Font.Bold = $true

$ws.Range("D10").Value = "В рамках UserService"
$ws.Range("D11").Value = "UserLegalService+Repository"

# B19 (admin.site.travel label) becomes bold, matching its sibling B18.
$ws.Range("B19").Font.Bold = $true

# Column C got wider to fit the new text.
$ws.Columns.Item(3).ColumnWidth = 20.6666667

# Scroll/selection bookkeeping: the sheet view now shows row 26 selected
# (full-row selection), with the window scrolled down.
$ws.Range("A26:XFD26").Select()

# Restore the originally-active sheet ("Классы", 3rd tab) so the workbook-level
# active tab is unchanged by the above selection on "Этапность".
$wb.Worksheets.Item(3).Activate()
